$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output $ws.GetType().Name
$win = $excel.ActiveWindow
Write-Output $win.GetType().Name
